$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data update: "ProQuest" hit count corrected 190 -> 180 (G1 SUM(D:D) recalculates automatically) ---
$ws.Range("D2").Value = 180

# --- Column E was auto best-fit to the date column; widen it explicitly (drops the best-fit flag) ---
$ws.Columns.Item(5).ColumnWidth = 14.25

# --- B7 (exported "mit Benefits Advanced" search string) now wraps, with an explicit 126pt row height ---
$ws.Range("B7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 126

# --- Leave the cursor on the newly edited cell, matching the saved view state ---
$ws.Range("B7").Select()
